$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "defaultA2" -> "asfgqw"
$ws.Range("B2").Value = "asfgqw"

# A3: numeric 2 -> text "2" (store as text, not number)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2"
$ws.Range("A3").ClearFormats()

# C3: "defaultB3" -> "asfagsad"
$ws.Range("C3").Value = "asfagsad"

# D3: "defaultB4" -> "hsdfqvxz"
$ws.Range("D3").Value = "hsdfqvxz"

# H4: "defaultC8" -> "fdsgcxvwsg"
$ws.Range("H4").Value = "fdsgcxvwsg"

# E5: "defaultD5" -> "sfsdfsdff"
$ws.Range("E5").Value = "sfsdfsdff"

# H5: "defaultD8" -> "defaultD8fd"
$ws.Range("H5").Value = "defaultD8fd"

# D6: "defaultE4" -> "фыыфвфы" (Cyrillic)
$ws.Range("D6").Value = "фыыфвфы"
